$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraph: remove the "_GoBack" bookmark that splits the
#    word "Instructions" into " I" / "nstructions" runs, while keeping
#    the visible text the same ("In-Class Task 8 Instructions: ...").
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1).Range
$titlePara.Find.Execute("8 Instructions:", $false, $false, $false, $false, $false, $true, 1, $false, "8 Instructions:", 1) | Out-Null

# ------------------------------------------------------------------
# 2) Second paragraph: drop the proofErr gramStart/gramEnd markers
#    around "application" and add a comma after it
#    ("...Invoice application so..." -> "...Invoice application, so...").
# ------------------------------------------------------------------
$secondPara = $d.Paragraphs(2).Range
$secondPara.Find.Execute("Invoice application so", $false, $false, $false, $false, $false, $true, 1, $false, "Invoice application, so", 1) | Out-Null

# ------------------------------------------------------------------
# 3) Bold the "Code a function..." list item.
# ------------------------------------------------------------------
$codeRange = $d.Content
$found = $codeRange.Find.Execute("Code a function that formats the Date object", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $codePara = $codeRange.Paragraphs(1).Range
    $codePara.Bold = 1
    $codePara.BoldBi = 1
}
